$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vtab = [char]11

$t.Cell(1, 1).Range.Text = "42 x 20" + $vtab + "  2    0" + $vtab + "  ----" + $vtab + "4|    |" + $vtab + "2|    |"
$t.Cell(1, 2).Range.Text = "56 x 97" + $vtab + "  9    7" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "6|    |"
$t.Cell(1, 3).Range.Text = "85 x 19" + $vtab + "  1    9" + $vtab + "  ----" + $vtab + "8|    |" + $vtab + "5|    |"
$t.Cell(2, 1).Range.Text = "73 x 88" + $vtab + "  8    8" + $vtab + "  ----" + $vtab + "7|    |" + $vtab + "3|    |"
$t.Cell(2, 2).Range.Text = "10 x 30" + $vtab + "  3    0" + $vtab + "  ----" + $vtab + "1|    |" + $vtab + "0|    |"
$t.Cell(2, 3).Range.Text = "64 x 30" + $vtab + "  3    0" + $vtab + "  ----" + $vtab + "6|    |" + $vtab + "4|    |"
$t.Cell(3, 1).Range.Text = "33 x 41" + $vtab + "  4    1" + $vtab + "  ----" + $vtab + "3|    |" + $vtab + "3|    |"
$t.Cell(3, 2).Range.Text = "88 x 83" + $vtab + "  8    3" + $vtab + "  ----" + $vtab + "8|    |" + $vtab + "8|    |"
$t.Cell(3, 3).Range.Text = "96 x 52" + $vtab + "  5    2" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "6|    |"
$t.Cell(4, 1).Range.Text = "99 x 71" + $vtab + "  7    1" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "9|    |"
$t.Cell(4, 2).Range.Text = "14 x 52" + $vtab + "  5    2" + $vtab + "  ----" + $vtab + "1|    |" + $vtab + "4|    |"
$t.Cell(4, 3).Range.Text = "78 x 51" + $vtab + "  5    1" + $vtab + "  ----" + $vtab + "7|    |" + $vtab + "8|    |"
$t.Cell(5, 1).Range.Text = "88 x 27" + $vtab + "  2    7" + $vtab + "  ----" + $vtab + "8|    |" + $vtab + "8|    |"
$t.Cell(5, 2).Range.Text = "35 x 14" + $vtab + "  1    4" + $vtab + "  ----" + $vtab + "3|    |" + $vtab + "5|    |"
$t.Cell(5, 3).Range.Text = "67 x 53" + $vtab + "  5    3" + $vtab + "  ----" + $vtab + "6|    |" + $vtab + "7|    |"
